$d = $word.ActiveDocument

# --- Change 1: rewrite the paragraph starting "el de leden van de Eerste Kamer..." ---
# (continuation of the "Hoewel" heading) into multiple runs with updated wording.
$target1 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -match "rechtstreeks op hen") {
        $target1 = $cand
        break
    }
}

$full = $target1.Range
$pStart = $full.Start
$pEnd = $full.End - 1  # exclude the paragraph mark

$pieces = @(
    "D",
    "e leden van de Eerste Kamer ook worden",
    " ",
    "gekozen",
    ", stemmen we niet rechtstreeks op hen. De leden van de Eerste Kamer worden gekozen door de",
    " andere",
    " leden van de ",
    "staten van de provincies",
    ". Na de Provinciale Statenverkiezingen vormen de nieuw gekozen Statenleden de kiesgroepen die op hun beurt de leden van de Eerste Kamer kiezen."
)

# Replace the whole paragraph's text with the first piece (becomes run 1)
$target = $d.Range($pStart, $pEnd)
$target.Text = $pieces[0]
$pos = $pStart + $pieces[0].Length

# Insert each subsequent piece and force it into its own run boundary.
for ($i = 1; $i -lt $pieces.Count; $i++) {
    $piece = $pieces[$i]
    $insRange = $d.Range($pos, $pos)
    $insRange.InsertAfter($piece)
    $newRange = $d.Range($pos, $pos + $piece.Length)
    $newRange.Bold = 1
    $newRange.Bold = 0
    $pos = $pos + $piece.Length
}

# --- Change 2: add a new empty paragraph at the very end of the document body,
# right before the final section break. ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$endPos = $lastPara.Range.End
$markRange = $d.Range($endPos, $endPos)
$markRange.InsertBefore("`r")

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
